$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 6, shifting existing rows 6-43 down to 7-44
$ws.Rows.Item(6).Insert()

# Populate the newly inserted row 6 with the new price-report entry
$ws.Range("A6").Value = 5
$ws.Range("B6").Value = 'Macroferia Regional de Talca'
$ws.Range("C6").Value = 'Maule'
$ws.Range("D6").Value = 45022
$ws.Range("E6").Value = 7
$ws.Range("F6").Value = 'Fruta'
$ws.Range("G6").Value = 100107
$ws.Range("H6").Value = 'Otros'
$ws.Range("I6").Value = 100107011
$ws.Range("J6").Value = 'Tuna'
$ws.Range("K6").Value = 'Sin especificar'
$ws.Range("L6").Value = 'Especial'
$ws.Range("M6").Value = 200
$ws.Range("N6").Value = 18000
$ws.Range("O6").Value = 18000
$ws.Range("P6").Value = 18000
$ws.Range("Q6").Value = '$/caja 18 kilos'
$ws.Range("R6").Value = 'Provincia de Melipilla'
$ws.Range("S6").Value = 1000
$ws.Range("T6").Value = 18
